$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header banner (rows 1-2) ---
# E1 needs the "title" banner style (same as A1/A2) before getting its new value.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "BILL OF MATERIALS"

# D1 drops from the bold/size-20-on-dark style down to the plain dark-fill style (like F1).
$ws.Range("F1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Row 2: board name text updates; date + revision swap columns.
$ws.Range("A2").Value = "BOARD: Z- SUN SENSE BOARD"

$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "REVISION: A"

$ws.Range("F2").Value = "Date: 2021.01.28"

$ws.Range("D2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Table header rename: Link -> Datasheet ---
$ws.Range("I3").Value = "Datasheet"

# --- Renumber the Item column (A4:A19) sequentially; old numbering skipped "2" ---
$items = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
for ($i = 0; $i -lt $items.Length; $i++) {
    $ws.Cells.Item(4 + $i, 1).Value = $items[$i]
}

# --- Drop the trailing, unused column J (dimension shrinks from J19 to I19) ---
$ws.Range("J:J").Delete()

# --- Match the saved selection ---
$ws.Range("A2").Select()
